$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 85 and 86 (id 6404135 / 6402965 odds refreshed + positions swapped)
# Row 85
$ws.Range("B85").Value = 6402965
$ws.Range("F85").Value = 'Maccabi Bnei Raina'
$ws.Range("G85").Value = 'Hapoel Hadera'
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = 1
$ws.Range("K85").Value = 2.625
$ws.Range("M85").Value = 2.4
$ws.Range("N85").Value = 2.6
$ws.Range("O85").Value = 3
$ws.Range("P85").Value = 2.7
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = 1.9
$ws.Range("S85").Value = 1.95
$ws.Range("T85").Value = 2
$ws.Range("U85").Value = 1.875
$ws.Range("V85").Value = 1.975
$ws.Range("X85").Value = 2
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Value = -0
$ws.Range("AB85").Value = 0
$ws.Range("AC85").Value = -0

# Row 86
$ws.Range("B86").Value = 6404135
$ws.Range("F86").Value = 'Hapoel Kiryat Shmona'
$ws.Range("G86").Value = 'Hapoel Bnei Sakhnin'
$ws.Range("H86").Value = 2
$ws.Range("I86").Value = 2
$ws.Range("K86").Value = 2.7
$ws.Range("M86").Value = 2.375
$ws.Range("N86").Value = 2
$ws.Range("O86").Value = 3.4
$ws.Range("P86").Value = 3.25
$ws.Range("Q86").Value = -0.25
$ws.Range("R86").Value = 1.8
$ws.Range("S86").Value = 2.05
$ws.Range("T86").Value = 2.25
$ws.Range("U86").Value = 1.8
$ws.Range("V86").Value = 2.05
$ws.Range("X86").Value = 2.4
$ws.Range("Z86").Value = -0.5
$ws.Range("AA86").Value = 0.5249999999999999
$ws.Range("AB86").Value = 0.8
$ws.Range("AC86").Value = -1

# Row 174
$ws.Range("B174").Value = 7511181
$ws.Range("F174").Value = 'Hapoel Hadera'
$ws.Range("G174").Value = 'Maccabi Bnei Raina'
$ws.Range("H174").Value = 0
$ws.Range("I174").Value = 2
$ws.Range("J174").Value = 'A'
$ws.Range("K174").Value = 2.7
$ws.Range("L174").Value = 3.3
$ws.Range("M174").Value = 2.4
$ws.Range("N174").Value = 2.45
$ws.Range("O174").Value = 3.25
$ws.Range("P174").Value = 2.7
$ws.Range("R174").Value = 1.825
$ws.Range("S174").Value = 2.025
$ws.Range("T174").Value = 2.25
$ws.Range("U174").Value = 1.875
$ws.Range("V174").Value = 1.975
$ws.Range("W174").Value = -1
$ws.Range("Y174").Value = 1.7
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 1.025
$ws.Range("AB174").Value = -0.5
$ws.Range("AC174").Value = 0.4875

# Row 175
$ws.Range("B175").Value = 7511180
$ws.Range("F175").Value = 'Maccabi Petach Tikva'
$ws.Range("G175").Value = 'Maccabi Netanya'
$ws.Range("H175").Value = 1
$ws.Range("I175").Value = 0
$ws.Range("J175").Value = 'H'
$ws.Range("K175").Value = 2.5
$ws.Range("L175").Value = 3.2
$ws.Range("M175").Value = 2.5
$ws.Range("N175").Value = 2.5
$ws.Range("O175").Value = 3.2
$ws.Range("P175").Value = 2.45
$ws.Range("R175").Value = 1.95
$ws.Range("S175").Value = 1.9
$ws.Range("T175").Value = 2.5
$ws.Range("U175").Value = 2
$ws.Range("V175").Value = 1.85
$ws.Range("W175").Value = 1.5
$ws.Range("Y175").Value = -1
$ws.Range("Z175").Value = 0.95
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = -1
$ws.Range("AC175").Value = 0.8500000000000001

# Row 202
$ws.Range("B202").Value = 7542719
$ws.Range("F202").Value = 'Hapoel Haifa'
$ws.Range("G202").Value = 'Maccabi Netanya'
$ws.Range("H202").Value = 2
$ws.Range("J202").Value = 'H'
$ws.Range("K202").Value = 2.6
$ws.Range("L202").Value = 3.1
$ws.Range("M202").Value = 2.6
$ws.Range("N202").Value = 2.9
$ws.Range("O202").Value = 3.2
$ws.Range("P202").Value = 2.3
$ws.Range("Q202").Value = 0.25
$ws.Range("R202").Value = 1.8
$ws.Range("S202").Value = 2.05
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 2
$ws.Range("V202").Value = 1.85
$ws.Range("W202").Value = 1.9
$ws.Range("X202").Value = -1
$ws.Range("Z202").Value = 0.8
$ws.Range("AA202").Value = -1
$ws.Range("AB202").Value = 1
$ws.Range("AC202").Value = -1

# Row 203
$ws.Range("B203").Value = 7542640
$ws.Range("F203").Value = 'MS Ashdod'
$ws.Range("G203").Value = 'Hapoel Bnei Sakhnin'
$ws.Range("H203").Value = 0
$ws.Range("J203").Value = 'A'
$ws.Range("K203").Value = 2.05
$ws.Range("L203").Value = 3.2
$ws.Range("M203").Value = 3.5
$ws.Range("N203").Value = 2.15
$ws.Range("O203").Value = 3.1
$ws.Range("P203").Value = 3.2
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 1.925
$ws.Range("S203").Value = 1.925
$ws.Range("T203").Value = 2.25
$ws.Range("U203").Value = 1.9
$ws.Range("V203").Value = 1.95
$ws.Range("W203").Value = -1
$ws.Range("Y203").Value = 2.2
$ws.Range("Z203").Value = -1
$ws.Range("AA203").Value = 0.925
$ws.Range("AB203").Value = -1
$ws.Range("AC203").Value = 0.95

# Row 204
$ws.Range("B204").Value = 7542639
$ws.Range("F204").Value = 'Maccabi Bnei Raina'
$ws.Range("G204").Value = 'Hapoel Jerusalem FC'
$ws.Range("H204").Value = 1
$ws.Range("J204").Value = 'D'
$ws.Range("K204").Value = 2.5
$ws.Range("L204").Value = 3
$ws.Range("M204").Value = 2.75
$ws.Range("N204").Value = 2.7
$ws.Range("O204").Value = 2.8
$ws.Range("P204").Value = 2.75
$ws.Range("Q204").Value = 0
$ws.Range("T204").Value = 2
$ws.Range("U204").Value = 2.1
$ws.Range("V204").Value = 1.775
$ws.Range("X204").Value = 1.8
$ws.Range("Y204").Value = -1
$ws.Range("Z204").Value = 0
$ws.Range("AA204").Value = -0
$ws.Range("AB204").Value = 0
$ws.Range("AC204").Value = -0

# Row 216
$ws.Range("B216").Value = 7542726
$ws.Range("F216").Value = 'Hapoel Hadera'
$ws.Range("G216").Value = 'Maccabi Netanya'
$ws.Range("H216").Value = 1
$ws.Range("I216").Value = 4
$ws.Range("K216").Value = 3.3
$ws.Range("L216").Value = 3.5
$ws.Range("M216").Value = 2
$ws.Range("N216").Value = 4.333
$ws.Range("O216").Value = 3.6
$ws.Range("P216").Value = 1.7
$ws.Range("Q216").Value = 0.75
$ws.Range("R216").Value = 1.9
$ws.Range("S216").Value = 1.95
$ws.Range("T216").Value = 2.5
$ws.Range("U216").Value = 2
$ws.Range("V216").Value = 1.85
$ws.Range("Y216").Value = 0.7
$ws.Range("AA216").Value = 0.95
$ws.Range("AB216").Value = 1
$ws.Range("AC216").Value = -1

# Row 217
$ws.Range("B217").Value = 7542727
$ws.Range("F217").Value = 'Maccabi Bnei Raina'
$ws.Range("G217").Value = 'Hapoel Bnei Sakhnin'
$ws.Range("H217").Value = 0
$ws.Range("I217").Value = 1
$ws.Range("K217").Value = 2.1
$ws.Range("L217").Value = 3.1
$ws.Range("M217").Value = 3.6
$ws.Range("N217").Value = 2.45
$ws.Range("O217").Value = 3
$ws.Range("P217").Value = 3
$ws.Range("Q217").Value = -0.25
$ws.Range("R217").Value = 2.075
$ws.Range("S217").Value = 1.725
$ws.Range("T217").Value = 2.25
$ws.Range("U217").Value = 2.05
$ws.Range("V217").Value = 1.8
$ws.Range("Y217").Value = 2
$ws.Range("AA217").Value = 0.7250000000000001
$ws.Range("AB217").Value = -1
$ws.Range("AC217").Value = 0.8

# Row 279
$ws.Range("B279").Value = 6799970
$ws.Range("E279").Value = 45340.64583333334
$ws.Range("F279").Value = 'Maccabi Haifa'
$ws.Range("G279").Value = 'Hapoel Haifa'
$ws.Range("K279").Value = 1.363
$ws.Range("L279").Value = 4.6
$ws.Range("M279").Value = 7.5
$ws.Range("N279").Value = 1.363
$ws.Range("P279").Value = 8
$ws.Range("R279").Value = 1.925
$ws.Range("S279").Value = 1.925
$ws.Range("U279").Value = 1.925
$ws.Range("V279").Value = 1.925

# Row 280
$ws.Range("B280").Value = 6799968
$ws.Range("E280").Value = 45341.625
$ws.Range("F280").Value = 'Maccabi Netanya'
$ws.Range("G280").Value = 'Maccabi Tel Aviv'
$ws.Range("K280").Value = 6
$ws.Range("L280").Value = 4
$ws.Range("M280").Value = 1.5
$ws.Range("N280").Value = 6.5
$ws.Range("O280").Value = 4
$ws.Range("P280").Value = 1.5
$ws.Range("Q280").Value = 1
$ws.Range("R280").Value = 1.975
$ws.Range("S280").Value = 1.875
$ws.Range("T280").Value = 2.75

# Remove rows that correspond to cancelled/dropped fixtures (old rows 279-282,
# ids 6799967, 6799972, 6799971, 6799969). The remaining two fixtures (ids
# 6799970 and 6799968) were moved up into rows 279/280 above with refreshed
# odds, so delete the now-duplicate trailing rows 281-284.
$ws.Range("A281:A284").EntireRow.Delete() | Out-Null
